$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 267; existing rows 267-357 shift down to 268-358.
$ws.Rows.Item(267).Insert()

# Populate the newly inserted row 267 with the new weekly record.
$ws.Cells.Item(267, 1).Value = 11
$ws.Cells.Item(267, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(267, 3).Value = "Bíobío"
$ws.Cells.Item(267, 4).Value = 44876
$ws.Cells.Item(267, 5).Value = 8
$ws.Cells.Item(267, 6).Value = 100114001
$ws.Cells.Item(267, 7).Value = "Papa"
$ws.Cells.Item(267, 8).Value = "Asterix"
$ws.Cells.Item(267, 9).Value = "1a nueva(o)"
$ws.Cells.Item(267, 10).Value = 220
$ws.Cells.Item(267, 11).Value = 11000
$ws.Cells.Item(267, 12).Value = 12000
$ws.Cells.Item(267, 13).Value = 11545
$ws.Cells.Item(267, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(267, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(267, 16).Value = 462
$ws.Cells.Item(267, 17).Value = 25
$ws.Cells.Item(267, 18).Value = "Hortaliza"
